$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 45, shifting rows 45:90 down to 46:91
$ws.Rows.Item(45).EntireRow.Insert()

# Populate the new row 45 with a new price observation (same market/product
# metadata as every other row in this sheet; only the date/volume/price
# columns differ per record).
$ws.Range("A45").Value = 10
$ws.Range("B45").Value = "Vega Modelo de Temuco"
$ws.Range("C45").Value = "La Araucanía"
$ws.Range("D45").Value = 44778
$ws.Range("E45").Value = 9
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100108
$ws.Range("H45").Value = "Tropicales y subtropicales"
$ws.Range("I45").Value = 100108007
$ws.Range("J45").Value = "Coco"
$ws.Range("K45").Value = "Sin especificar"
$ws.Range("L45").Value = "Primera"
$ws.Range("M45").Value = 15
$ws.Range("N45").Value = 30000
$ws.Range("O45").Value = 30000
$ws.Range("P45").Value = 30000
$ws.Range("Q45").Value = "$/malla 20 unidades"
$ws.Range("R45").Value = "Perú"
$ws.Range("S45").Value = 1500
$ws.Range("T45").Value = 20

Write-Output "done"
